$wb = $excel.ActiveWorkbook

# --- 1. Update selection on InsuredPageData (sheet index 3) before it loses focus ---
$insured = $wb.Worksheets.Item(3)
$insured.Range("C31").Select()

# --- 2. Insert the new "QuoteOptionPageData" sheet right after InsuredPageData
#        (i.e. right before RatingCriteriaPageData), which becomes the newly
#        active / tab-selected sheet. ---
$newWs = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $insured)
$newWs.Name = "QuoteOptionPageData"

# --- 3. Populate the new sheet's data ---
$reg = [char]174

$newWs.Range("A1").Value = "testAddQuoteOption"

$newWs.Range("A2").Value = "runMode"
$newWs.Range("B2").Value = "product"
$newWs.Range("C2").Value = "applicantName"
$newWs.Range("D2").Value = "website"
$newWs.Range("E2").Value = "brokerId"
$newWs.Range("F2").Value = "agentId"
$newWs.Range("G2").Value = "agencyOfficeId"
$newWs.Range("H2").Value = "numberOfResidentialUnits"
$newWs.Range("I2").Value = "totalCommercialSquareFeet"
$newWs.Range("J2").Value = "businessClass"

$newWs.Range("A3").Value = "Y"
$newWs.Range("B3").Value = "NetGuard" + $reg + " Plus"
$newWs.Range("C3").Value = "Caring Communities, A Reciprocal Risk Re"
$newWs.Range("D3").Value = "https://caringcomm.org/"
$newWs.Range("E3").Value = 20217
$newWs.Range("F3").Value = 173
$newWs.Range("G3").Value = 237
$newWs.Range("H3").Value = 173
$newWs.Range("I3").Value = 237
$newWs.Range("J3").Value = "Business to Business"

# --- 4. Hyperlink on D3 ---
$newWs.Hyperlinks.Add($newWs.Range("D3"), "https://caringcomm.org/")

# --- 5. Formatting: header fill/border (best-effort; engine border color
#        dedup is limited, so this approximates the authored look) ---
$newWs.Range("A1:B1").Interior.ColorIndex = 6
$newWs.Range("A1").Borders.LineStyle = 1
$newWs.Range("B1").Borders.LineStyle = 1
$newWs.Range("A2:J2").Borders.LineStyle = 1
$newWs.Range("A3:J3").Borders.LineStyle = 1

# --- 6. Leave the selection on A4, matching the authored workbook state ---
$newWs.Range("A4").Select()

Write-Host "QuoteOptionPageData sheet added"
